# The "Farbe:text*; Kaufdatum:date" attribute row is rewritten to use
# ":pflicht" instead of a trailing "*" to mark the field as mandatory, and
# likewise for the "Typ:select(...)" attribute row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Farbe:text:pflicht; Kaufdatum:date"
$ws.Range("D4").Value = "Typ:select(Auto,Fahrrad,Roller):pflicht"

# Reflect the author's final on-screen selection/scroll state: view was
# scrolled right one column (topLeftCell C1) with D2 the active cell.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("D2").Select()
